$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.009607864071965899, 0.02547733431879777, 0.0111111111111111, 0, 0.005191550304332257, 0.01352129911543065, 0.01004851004851003, 0.006300553426990201, 0.006098191214470299, 0.01675824175824175),
    @(0.01861105608723967, 0.01049943883277217, 0.01363740497761117, 0.02104375460636393, 0.02275406833925, 0.02121957727909081, 0.01429270583243773, 0.005339105339105334, 0.008233371689181033, 0.0233184872629388),
    @(0.02223648218185563, 0.01223717733221533, 0.0176162416708393, 0.02263095373485397, 0.02792940243643537, 0.0206926638517629, 0.01884081015136253, 0.0064031950386499, 0.006074723297124736, 0.02131717203729095),
    @(0.03242169269328047, 0.01994268055928894, 0.02538135047134923, 0.02514440029998933, 0.0232964879515586, 0.0254706311100013, 0.0165216791157632, 0.005387159205103467, 0.004505383996059827, 0.02500390686044695),
    @(0.02106224659716753, 0.0179622236038704, 0.015284797178501, 0.01427922215176497, 0.0151698288310754, 0.02032942032942034, 0.0131578947368421, 0.0060309351497091, 0.005847953216374266, 0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}
